# Updates the cryptos price/volume table with newly scraped values.
# Cells are stored as plain (non-numeric) text in the sheet, so each write
# forces a Text number format before assigning the value (otherwise Excel's
# COM layer auto-coerces number-looking strings like "1.00" into the number 1)
# and then clears the format again so no stray number-format style is left
# behind on the cell (matching the original, unstyled cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "58.814.38"
Set-TextValue $ws.Range("E2") "  +2.59%  "
Set-TextValue $ws.Range("D3") "2.543.65"
Set-TextValue $ws.Range("E3") "  +4.96%  "
Set-TextValue $ws.Range("E4") "  -0.13%  "
Set-TextValue $ws.Range("D5") "527.42"
Set-TextValue $ws.Range("E5") "  +3.08%  "
Set-TextValue $ws.Range("D6") "134.46"
Set-TextValue $ws.Range("E6") "  +2.63%  "
Set-TextValue $ws.Range("D7") "1.00"
Set-TextValue $ws.Range("E7") "  +0.08%  "
Set-TextValue $ws.Range("D8") "0.566"
Set-TextValue $ws.Range("E8") "  +2.88%  "
Set-TextValue $ws.Range("D9") "2.541.95"
Set-TextValue $ws.Range("E9") "  +4.85%  "
Set-TextValue $ws.Range("D10") "0.0991"
Set-TextValue $ws.Range("E10") "  +3.61%  "
Set-TextValue $ws.Range("E11") "  -0.86%  "
Set-TextValue $ws.Range("E12") "  -0.08%  "
Set-TextValue $ws.Range("E13") "  +1.80%  "
Set-TextValue $ws.Range("D14") "2.995.64"
Set-TextValue $ws.Range("E14") "  +4.89%  "
Set-TextValue $ws.Range("D15") "58.848.63"
Set-TextValue $ws.Range("E15") "  +2.75%  "
Set-TextValue $ws.Range("D16") "22.43"
Set-TextValue $ws.Range("E16") "  +4.29%  "
Set-TextValue $ws.Range("E17") "  +3.60%  "
Set-TextValue $ws.Range("D18") "2.541.11"
Set-TextValue $ws.Range("E18") "  +4.32%  "
Set-TextValue $ws.Range("D19") "10.74"
Set-TextValue $ws.Range("E19") "  +3.78%  "
Set-TextValue $ws.Range("D20") "324.15"
Set-TextValue $ws.Range("E20") "  +3.18%  "
Set-TextValue $ws.Range("E21") "  +3.10%  "
Set-TextValue $ws.Range("E22") "  +10.23%  "
Set-TextValue $ws.Range("E23") "  +0.18%  "
Set-TextValue $ws.Range("D24") "65.13"
Set-TextValue $ws.Range("E24") "  +2.22%  "
Set-TextValue $ws.Range("D25") "0.411"
Set-TextValue $ws.Range("E25") "  +2.54%  "
Set-TextValue $ws.Range("D26") "0.999"
Set-TextValue $ws.Range("E26") "  -0.29%  "
Set-TextValue $ws.Range("E27") "  +1.72%  "
Set-TextValue $ws.Range("D28") "7.48"
Set-TextValue $ws.Range("E28") "  +4.22%  "
Set-TextValue $ws.Range("E29") "  +4.95%  "
Set-TextValue $ws.Range("E30") "  +7.12%  "
Set-TextValue $ws.Range("E31") "  +4.59%  "
Set-TextValue $ws.Range("D32") "168.88"
Set-TextValue $ws.Range("E32") "  -0.28%  "
Set-TextValue $ws.Range("D33") "6.37"
Set-TextValue $ws.Range("E33") "  +2.87%  "
Set-TextValue $ws.Range("D34") "0.998"
Set-TextValue $ws.Range("E34") "  -0.06%  "
Set-TextValue $ws.Range("D35") "0.996"
Set-TextValue $ws.Range("E35") "  -0.29%  "
Set-TextValue $ws.Range("D36") "18.29"
Set-TextValue $ws.Range("E36") "  +3.53%  "
Set-TextValue $ws.Range("D37") "1.27"
Set-TextValue $ws.Range("E37") "  -0.89%  "
Set-TextValue $ws.Range("E38") "  +3.42%  "
Set-TextValue $ws.Range("D39") "1.52"
Set-TextValue $ws.Range("E39") "  +5.89%  "
Set-TextValue $ws.Range("D40") "36.84"
Set-TextValue $ws.Range("E40") "  +1.75%  "
Set-TextValue $ws.Range("E41") "  +1.74%  "
Set-TextValue $ws.Range("D42") "281.41"
Set-TextValue $ws.Range("E42") "  +5.89%  "
Set-TextValue $ws.Range("D43") "3.48"
Set-TextValue $ws.Range("E43") "  +4.03%  "
Set-TextValue $ws.Range("D44") "133.69"
Set-TextValue $ws.Range("E44") "  +10.01%  "
Set-TextValue $ws.Range("E45") "  +3.86%  "
Set-TextValue $ws.Range("E46") "  +4.10%  "
Set-TextValue $ws.Range("E47") "  +2.77%  "
Set-TextValue $ws.Range("D48") "0.0507"
Set-TextValue $ws.Range("E48") "  +5.76%  "
Set-TextValue $ws.Range("D49") "17.87"
Set-TextValue $ws.Range("E49") "  +4.60%  "
Set-TextValue $ws.Range("E50") "  +4.02%  "
Set-TextValue $ws.Range("D51") "17.20"
Set-TextValue $ws.Range("E51") "  +4.25%  "
